$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete old sample data rows (2-5) entirely so the sheet only keeps the header row
$ws.Range("A2:K5").Delete()
$ws.Rows.Ungroup()

# 2. Update / add header row values (A1:P1) in the new column order
$headers = @(
    "stt",
    "ma_so_bhxh",
    "ho_ten_hoc_sinh",
    "ngay_sinh",
    "gioi_tinh",
    "dia_chi",
    "ngay_het_han_bhyt",
    "ngay_het_han_bhtn",
    "lop_hoc",
    "sdt_lienhe",
    "so_dinh_danh",
    "noi_kham_bhyt",
    "ten_cha_me",
    "doi_tuong_dong",
    "ghi_chu",
    "ma_truong"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $ws.Range("A1:P1")

# 3. Make sure every header cell shares the same font / wrap formatting
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 10
$headerRange.WrapText = $true
$headerRange.VerticalAlignment = -4107

# 4. Adjust column widths
$ws.Columns.Item(7).ColumnWidth = 20.0                      # G  ngay_het_han_bhyt
$ws.Columns.Item(8).ColumnWidth = 19.0                       # H  ngay_het_han_bhtn
$ws.Columns.Item(12).ColumnWidth = 19.666666666666668         # L  noi_kham_bhyt
$ws.Columns.Item(13).ColumnWidth = 14.0                       # M  ten_cha_me
$ws.Columns.Item(14).ColumnWidth = 13.0                       # N  doi_tuong_dong
$ws.Columns.Item(15).ColumnWidth = 9.833333333333332          # O  ghi_chu
$ws.Columns.Item(16).ColumnWidth = 11.166666666666666         # P  ma_truong

# 5. Apply a medium grey (#CCCCCC) border around all header cells
$headerRange.Borders.Weight = -4138
$headerRange.Borders.Color = 13421772

# 6. Update selection to match the saved view
$ws.Range("F16").Select() | Out-Null

Write-Output "done"
